$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.6268010600226682
$ws.Range("C4").Value = 0.634
$ws.Range("D4").Value = 0.6281201169652636
$ws.Range("E4").Value = 0.6234999999999999
$ws.Range("F4").Value = 0.5087558688532624
$ws.Range("G4").Value = 0.515
$ws.Range("H4").Value = 0.5044277824880226
$ws.Range("I4").Value = 0.505
$ws.Range("J4").Value = 0.6518307222623771
$ws.Range("K4").Value = 0.6719999999999999
$ws.Range("L4").Value = 0.6365815140379242
$ws.Range("M4").Value = 0.6425000000000001

$ws.Range("B5").Value = 0.6716795498415088
$ws.Range("C5").Value = 0.8400000000000001
$ws.Range("D5").Value = 0.5681587307412322
$ws.Range("E5").Value = 0.591
$ws.Range("F5").Value = 0.6682751163541797
$ws.Range("G5").Value = 0.9700000000000001
$ws.Range("H5").Value = 0.5102262724347579
$ws.Range("I5").Value = 0.518
$ws.Range("J5").Value = 0.6475559102713125
$ws.Range("K5").Value = 0.8280000000000001
$ws.Range("L5").Value = 0.5375566682583701
$ws.Range("M5").Value = 0.5545

$ws.Range("B6").Value = 0.6362132927802666
$ws.Range("C6").Value = 0.649
$ws.Range("D6").Value = 0.6306867986216904
$ws.Range("E6").Value = 0.6315
$ws.Range("F6").Value = 0.4925158945170987
$ws.Range("G6").Value = 0.4869999999999999
$ws.Range("H6").Value = 0.5006483325187048
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.6438320996789805
$ws.Range("K6").Value = 0.6240000000000001
$ws.Range("L6").Value = 0.6769027801196138
$ws.Range("M6").Value = 0.659
